# Add a "user_id" column (Telegram user id) to the users sheet, between
# "id" and "api_id", and populate it for the three existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (api_id) -> new column B = user_id
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "user_id"
$ws.Range("B2").Value = 784093829
$ws.Range("B3").Value = 1442776649
$ws.Range("B4").Value = 1396995011

# Match the column width from the diff for the new column B (Excel
# auto-fits to the longest value, same as double-clicking the column
# border). ColumnWidth uses "characters" units which Excel stores with a
# small fixed padding added on save, so back the padding out here to land
# on the serialized width of 11 recorded in the target file.
$ws.Columns.Item(2).ColumnWidth = 10.1666666666667

# Update selection to match the committed workbook state.
$ws.Range("C9").Select()
